$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '63.725.27'
$ws.Cells.Item(2, 5).Value = '  +1.64%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.474.33'
$ws.Cells.Item(3, 5).Value = '  +1.35%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '575.45'
$ws.Cells.Item(5, 5).Value = '  +1.47%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '148.64'
$ws.Cells.Item(6, 5).Value = '  +2.18%  '
$ws.Cells.Item(7, 5).Value = '  +0.12%  '
$ws.Cells.Item(8, 5).Value = '  +1.80%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.474.50'
$ws.Cells.Item(9, 5).Value = '  +1.31%  '
$ws.Cells.Item(10, 5).Value = '  +0.93%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.156'
$ws.Cells.Item(11, 5).Value = '  +1.15%  '
$ws.Cells.Item(12, 5).Value = '  +1.16%  '
$ws.Cells.Item(13, 5).Value = '  +1.50%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '27.30'
$ws.Cells.Item(14, 5).Value = '  +1.16%  '
$ws.Cells.Item(15, 5).Value = '  -0.42%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '2.924.53'
$ws.Cells.Item(16, 5).Value = '  +1.51%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.596.25'
$ws.Cells.Item(17, 5).Value = '  +1.68%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.474.55'
$ws.Cells.Item(18, 5).Value = '  +1.81%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.46'
$ws.Cells.Item(19, 5).Value = '  +1.96%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.44'
$ws.Cells.Item(20, 5).Value = '  +7.31%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '331.61'
$ws.Cells.Item(21, 5).Value = '  +2.56%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.23'
$ws.Cells.Item(22, 5).Value = '  +1.55%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '2.11'
$ws.Cells.Item(23, 5).Value = '  +17.76%  '
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '66.01'
$ws.Cells.Item(25, 5).Value = '  -1.75%  '
$ws.Cells.Item(26, 2).Value = 'Aptos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.27'
$ws.Cells.Item(26, 5).Value = '  +6.85%  '
$ws.Cells.Item(27, 2).Value = 'Bittensor'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '631.69'
$ws.Cells.Item(27, 5).Value = '  +10.65%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0000105'
$ws.Cells.Item(28, 5).Value = '  +3.42%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.54'
$ws.Cells.Item(29, 5).Value = '  +6.04%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.605.16'
$ws.Cells.Item(30, 5).Value = '  +1.77%  '
$ws.Cells.Item(31, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.00'
$ws.Cells.Item(31, 5).Value = '  +0.11%  '
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '8.42'
$ws.Cells.Item(32, 5).Value = '  +0.38%  '
$ws.Cells.Item(33, 5).Value = '  -2.42%  '
$ws.Cells.Item(34, 5).Value = '  +2.85%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.28'
$ws.Cells.Item(35, 5).Value = '  +8.77%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.57'
$ws.Cells.Item(36, 5).Value = '  +1.15%  '
$ws.Cells.Item(37, 5).Value = '  +0.25%  '
$ws.Cells.Item(38, 5).Value = '  +0.20%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.56'
$ws.Cells.Item(39, 5).Value = '  +2.42%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '18.93'
$ws.Cells.Item(40, 5).Value = '  +0.68%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.73'
$ws.Cells.Item(41, 5).Value = '  +13.28%  '
$ws.Cells.Item(42, 5).Value = '  +0.15%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '145.97'
$ws.Cells.Item(43, 5).Value = '  -1.83%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '151.05'
$ws.Cells.Item(45, 5).Value = '  +1.93%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.79'
$ws.Cells.Item(46, 5).Value = '  +3.28%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '21.63'
$ws.Cells.Item(48, 5).Value = '  +1.09%  '
$ws.Cells.Item(49, 5).Value = '  +1.00%  '
$ws.Cells.Item(50, 5).Value = '  +2.79%  '
$ws.Cells.Item(51, 5).Value = '  -0.45%  '
